$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L21").Value = 1108.6
$ws1.Range("L29").Value = "1 de 27"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F21").Value = 3225.33
$ws2.Range("F29").Value = 4147.16

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D15").Value = 2862.56
$ws3.Range("E15").Value = -1224.56
$ws3.Range("F15").Value = 1.747594627594627
$ws3.Range("D19").Value = 10734.92
$ws3.Range("E19").Value = 12765.08093005039
$ws3.Range("F19").Value = 0.4568050883041809
